$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dim = $ws.UsedRange
$lastRow = $dim.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in @("G", "H")) {
        $cell = $ws.Range("$col$r")
        $val = $cell.Value2
        if ($val -ne $null -and $val -ne "N/A") {
            $parts = $val -split ":"
            $newParts = @()
            foreach ($p in $parts) {
                $newParts += "0x" + $p
            }
            $newVal = $newParts -join ":"
            $cell.Value = $newVal
        }
    }
}
